# ============================================================================
# Refresh the Leve price/profit snapshot columns (H:N) on each job sheet.
#
# This workbook has no formulas -- H:N are literal values written by an
# external market-data puller ("scheduled runner"). This script pokes in the
# refreshed numbers cell-by-cell, grouped by sheet and row:
#   H = currentAveragePrice      K = LevePriceNQ
#   I = currentAveragePriceNQ    L = LevePriceHQ
#   J = currentAveragePriceHQ    M = LeveProfitNQ   N = LeveProfitHQ
#
# One cell (LTW!N122) is cleared outright rather than reassigned: its new
# LevePriceHQ (L122) came back as 0 from the price pull, so -- matching the
# sheet's existing convention for "no valid HQ price" rows (see e.g. row 5) --
# the HQ profit cell is left blank instead of holding a stale number.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 52779.75   # H64: 48236.137 -> 52779.75
$ws.Cells.Item(64, 10).Value = 2868.5625   # J64: 2860.9443 -> 2868.5625
$ws.Cells.Item(64, 12).Value = 2868.5625   # L64: 2860.9443 -> 2868.5625
$ws.Cells.Item(64, 14).Value = -3364.5625   # N64: -3356.9443 -> -3364.5625
# Row 67
$ws.Cells.Item(67, 8).Value = 52779.75   # H67: 48236.137 -> 52779.75
$ws.Cells.Item(67, 10).Value = 2868.5625   # J67: 2860.9443 -> 2868.5625
$ws.Cells.Item(67, 12).Value = 2868.5625   # L67: 2860.9443 -> 2868.5625
$ws.Cells.Item(67, 14).Value = -4584.5625   # N67: -4576.9443 -> -4584.5625
# Row 76
$ws.Cells.Item(76, 8).Value = 3133.5833   # H76: 3140.3 -> 3133.5833
$ws.Cells.Item(76, 9).Value = 3000.75   # I76: 3001 -> 3000.75
$ws.Cells.Item(76, 11).Value = 3000.75   # K76: 3001 -> 3000.75
$ws.Cells.Item(76, 13).Value = -2685.75   # M76: -2686 -> -2685.75
# Row 79
$ws.Cells.Item(79, 8).Value = 3133.5833   # H79: 3140.3 -> 3133.5833
$ws.Cells.Item(79, 9).Value = 3000.75   # I79: 3001 -> 3000.75
$ws.Cells.Item(79, 11).Value = 3000.75   # K79: 3001 -> 3000.75
$ws.Cells.Item(79, 13).Value = -1908.75   # M79: -1909 -> -1908.75
# Row 137
$ws.Cells.Item(137, 8).Value = 1331.75   # H137: 1236.8485 -> 1331.75
$ws.Cells.Item(137, 9).Value = 824.3103599999999   # I137: 765.1142599999999 -> 824.3103599999999
$ws.Cells.Item(137, 10).Value = 1539.014   # J137: 1494.8281 -> 1539.014
$ws.Cells.Item(137, 11).Value = 2472.93108   # K137: 2295.34278 -> 2472.93108
$ws.Cells.Item(137, 12).Value = 4617.041999999999   # L137: 4484.4843 -> 4617.041999999999
$ws.Cells.Item(137, 13).Value = 77.06892000000016   # M137: 254.6572200000001 -> 77.06892000000016
$ws.Cells.Item(137, 14).Value = -9717.041999999999   # N137: -9584.4843 -> -9717.041999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 24793.195   # H32: 25109.553 -> 24793.195
$ws.Cells.Item(32, 9).Value = 25508.127   # I32: 25861.814 -> 25508.127
$ws.Cells.Item(32, 11).Value = 25508.127   # K32: 25861.814 -> 25508.127
$ws.Cells.Item(32, 13).Value = -25221.127   # M32: -25574.814 -> -25221.127
# Row 61
$ws.Cells.Item(61, 8).Value = 2196.4893   # H61: 2299.8865 -> 2196.4893
$ws.Cells.Item(61, 9).Value = 2119.9048   # I61: 2228.1025 -> 2119.9048
$ws.Cells.Item(61, 10).Value = 2839.8   # J61: 2859.8 -> 2839.8
$ws.Cells.Item(61, 11).Value = 2119.9048   # K61: 2228.1025 -> 2119.9048
$ws.Cells.Item(61, 12).Value = 2839.8   # L61: 2859.8 -> 2839.8
$ws.Cells.Item(61, 13).Value = -1907.9048   # M61: -2016.1025 -> -1907.9048
$ws.Cells.Item(61, 14).Value = -3263.8   # N61: -3283.8 -> -3263.8
# Row 63
$ws.Cells.Item(63, 8).Value = 3027   # H63: 2993.4736 -> 3027
$ws.Cells.Item(63, 9).Value = 2475.3845   # I63: 2469.2856 -> 2475.3845
$ws.Cells.Item(63, 11).Value = 2475.3845   # K63: 2469.2856 -> 2475.3845
$ws.Cells.Item(63, 13).Value = -1789.3845   # M63: -1783.2856 -> -1789.3845
# Row 66
$ws.Cells.Item(66, 8).Value = 3027   # H66: 2993.4736 -> 3027
$ws.Cells.Item(66, 9).Value = 2475.3845   # I66: 2469.2856 -> 2475.3845
$ws.Cells.Item(66, 11).Value = 12376.9225   # K66: 12346.428 -> 12376.9225
$ws.Cells.Item(66, 13).Value = -8944.922500000001   # M66: -8914.428 -> -8944.922500000001
# Row 132
$ws.Cells.Item(132, 8).Value = 8476467   # H132: 12197400 -> 8476467
$ws.Cells.Item(132, 9).Value = 9805564   # I132: 14707845 -> 9805564
$ws.Cells.Item(132, 10).Value = 3471.5   # J132: 3811.1428 -> 3471.5
$ws.Cells.Item(132, 11).Value = 29416692   # K132: 44123535 -> 29416692
$ws.Cells.Item(132, 12).Value = 10414.5   # L132: 11433.4284 -> 10414.5
$ws.Cells.Item(132, 13).Value = -29414162   # M132: -44121005 -> -29414162
$ws.Cells.Item(132, 14).Value = -15474.5   # N132: -16493.4284 -> -15474.5
# Row 136
$ws.Cells.Item(136, 8).Value = 2196.4893   # H136: 2299.8865 -> 2196.4893
$ws.Cells.Item(136, 9).Value = 2119.9048   # I136: 2228.1025 -> 2119.9048
$ws.Cells.Item(136, 10).Value = 2839.8   # J136: 2859.8 -> 2839.8
$ws.Cells.Item(136, 11).Value = 6359.714399999999   # K136: 6684.3075 -> 6359.714399999999
$ws.Cells.Item(136, 12).Value = 8519.400000000001   # L136: 8579.400000000001 -> 8519.400000000001
$ws.Cells.Item(136, 13).Value = -3809.714399999999   # M136: -4134.3075 -> -3809.714399999999
$ws.Cells.Item(136, 14).Value = -13619.4   # N136: -13679.4 -> -13619.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 2778.125   # H62: 2544.0312 -> 2778.125
$ws.Cells.Item(62, 9).Value = 2640.9092   # I62: 2445.16 -> 2640.9092
$ws.Cells.Item(62, 10).Value = 3080   # J62: 2897.1428 -> 3080
$ws.Cells.Item(62, 11).Value = 2640.9092   # K62: 2445.16 -> 2640.9092
$ws.Cells.Item(62, 12).Value = 3080   # L62: 2897.1428 -> 3080
$ws.Cells.Item(62, 13).Value = -2016.9092   # M62: -1821.16 -> -2016.9092
$ws.Cells.Item(62, 14).Value = -4328   # N62: -4145.1428 -> -4328
# Row 65
$ws.Cells.Item(65, 8).Value = 2778.125   # H65: 2544.0312 -> 2778.125
$ws.Cells.Item(65, 9).Value = 2640.9092   # I65: 2445.16 -> 2640.9092
$ws.Cells.Item(65, 10).Value = 3080   # J65: 2897.1428 -> 3080
$ws.Cells.Item(65, 11).Value = 13204.546   # K65: 12225.8 -> 13204.546
$ws.Cells.Item(65, 12).Value = 15400   # L65: 14485.714 -> 15400
$ws.Cells.Item(65, 13).Value = -10084.546   # M65: -9105.799999999999 -> -10084.546
$ws.Cells.Item(65, 14).Value = -21640   # N65: -20725.714 -> -21640
# Row 99
$ws.Cells.Item(99, 8).Value = 1617   # H99: 1585.2354 -> 1617
$ws.Cells.Item(99, 9).Value = 1674.8572   # I99: 1513.5 -> 1674.8572
$ws.Cells.Item(99, 10).Value = 1559.1428   # J99: 1687.7142 -> 1559.1428
$ws.Cells.Item(99, 11).Value = 1674.8572   # K99: 1513.5 -> 1674.8572
$ws.Cells.Item(99, 12).Value = 1559.1428   # L99: 1687.7142 -> 1559.1428
$ws.Cells.Item(99, 13).Value = -176.8571999999999   # M99: -15.5 -> -176.8571999999999
$ws.Cells.Item(99, 14).Value = -4555.1428   # N99: -4683.7142 -> -4555.1428
# Row 126
$ws.Cells.Item(126, 8).Value = 1617   # H126: 1585.2354 -> 1617
$ws.Cells.Item(126, 9).Value = 1674.8572   # I126: 1513.5 -> 1674.8572
$ws.Cells.Item(126, 10).Value = 1559.1428   # J126: 1687.7142 -> 1559.1428
$ws.Cells.Item(126, 11).Value = 5024.571599999999   # K126: 4540.5 -> 5024.571599999999
$ws.Cells.Item(126, 12).Value = 4677.428400000001   # L126: 5063.142599999999 -> 4677.428400000001
$ws.Cells.Item(126, 13).Value = -2554.571599999999   # M126: -2070.5 -> -2554.571599999999
$ws.Cells.Item(126, 14).Value = -9617.428400000001   # N126: -10003.1426 -> -9617.428400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Cells.Item(122, 8).Value = 3059.8572   # H122: 2780.0625 -> 3059.8572
$ws.Cells.Item(122, 9).Value = 700.225   # I122: 677.4186 -> 700.225
$ws.Cells.Item(122, 10).Value = 50252.5   # J122: 20862.8 -> 50252.5
$ws.Cells.Item(122, 11).Value = 6302.025000000001   # K122: 6096.7674 -> 6302.025000000001
$ws.Cells.Item(122, 12).Value = 452272.5   # L122: 187765.2 -> 452272.5
$ws.Cells.Item(122, 13).Value = -3852.025000000001   # M122: -3646.7674 -> -3852.025000000001
$ws.Cells.Item(122, 14).Value = -457172.5   # N122: -192665.2 -> -457172.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 267113.16   # H80: 221050 -> 267113.16
$ws.Cells.Item(80, 10).Value = 3243.75   # J80: 2912.5 -> 3243.75
$ws.Cells.Item(80, 12).Value = 3243.75   # L80: 2912.5 -> 3243.75
$ws.Cells.Item(80, 14).Value = -5239.75   # N80: -4908.5 -> -5239.75
# Row 83
$ws.Cells.Item(83, 8).Value = 267113.16   # H83: 221050 -> 267113.16
$ws.Cells.Item(83, 10).Value = 3243.75   # J83: 2912.5 -> 3243.75
$ws.Cells.Item(83, 12).Value = 16218.75   # L83: 14562.5 -> 16218.75
$ws.Cells.Item(83, 14).Value = -26202.75   # N83: -24546.5 -> -26202.75
# Row 126
$ws.Cells.Item(126, 8).Value = 2495.52   # H126: 2538.5833 -> 2495.52
$ws.Cells.Item(126, 9).Value = 2735.3333   # I126: 2941.2 -> 2735.3333
$ws.Cells.Item(126, 10).Value = 2274.1538   # J126: 2251 -> 2274.1538
$ws.Cells.Item(126, 11).Value = 8205.999899999999   # K126: 8823.599999999999 -> 8205.999899999999
$ws.Cells.Item(126, 12).Value = 6822.4614   # L126: 6753 -> 6822.4614
$ws.Cells.Item(126, 13).Value = -5735.999899999999   # M126: -6353.599999999999 -> -5735.999899999999
$ws.Cells.Item(126, 14).Value = -11762.4614   # N126: -11693 -> -11762.4614
# Row 132
$ws.Cells.Item(132, 8).Value = 2580.2856   # H132: 2614.54 -> 2580.2856
$ws.Cells.Item(132, 9).Value = 2177.9744   # I132: 2226.35 -> 2177.9744
$ws.Cells.Item(132, 10).Value = 4149.3   # J132: 4167.3 -> 4149.3
$ws.Cells.Item(132, 11).Value = 6533.9232   # K132: 6679.049999999999 -> 6533.9232
$ws.Cells.Item(132, 12).Value = 12447.9   # L132: 12501.9 -> 12447.9
$ws.Cells.Item(132, 13).Value = -4003.9232   # M132: -4149.049999999999 -> -4003.9232
$ws.Cells.Item(132, 14).Value = -17507.9   # N132: -17561.9 -> -17507.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 2250   # H7: 2199.7058 -> 2250
$ws.Cells.Item(7, 9).Value = 2000   # I7: 1807.5 -> 2000
$ws.Cells.Item(7, 10).Value = 2300   # J7: 3141 -> 2300
$ws.Cells.Item(7, 11).Value = 2000   # K7: 1807.5 -> 2000
$ws.Cells.Item(7, 12).Value = 2300   # L7: 3141 -> 2300
$ws.Cells.Item(7, 13).Value = -1888   # M7: -1695.5 -> -1888
$ws.Cells.Item(7, 14).Value = -2524   # N7: -3365 -> -2524
# Row 40
$ws.Cells.Item(40, 8).Value = 2699.3333   # H40: 2114.923 -> 2699.3333
$ws.Cells.Item(40, 9).Value = 2439.4   # I40: 1957.9166 -> 2439.4
$ws.Cells.Item(40, 11).Value = 2439.4   # K40: 1957.9166 -> 2439.4
$ws.Cells.Item(40, 13).Value = -2303.4   # M40: -1821.9166 -> -2303.4
# Row 82
$ws.Cells.Item(82, 8).Value = 7577070.5   # H82: 7578059 -> 7577070.5
$ws.Cells.Item(82, 9).Value = 1305.7142   # I82: 3111.5715 -> 1305.7142
$ws.Cells.Item(82, 10).Value = 20834660   # J82: 20834216 -> 20834660
$ws.Cells.Item(82, 11).Value = 1305.7142   # K82: 3111.5715 -> 1305.7142
$ws.Cells.Item(82, 12).Value = 20834660   # L82: 20834216 -> 20834660
$ws.Cells.Item(82, 13).Value = -944.7141999999999   # M82: -2750.5715 -> -944.7141999999999
$ws.Cells.Item(82, 14).Value = -20835382   # N82: -20834938 -> -20835382
# Row 85
$ws.Cells.Item(85, 8).Value = 7577070.5   # H85: 7578059 -> 7577070.5
$ws.Cells.Item(85, 9).Value = 1305.7142   # I85: 3111.5715 -> 1305.7142
$ws.Cells.Item(85, 10).Value = 20834660   # J85: 20834216 -> 20834660
$ws.Cells.Item(85, 11).Value = 1305.7142   # K85: 3111.5715 -> 1305.7142
$ws.Cells.Item(85, 12).Value = 20834660   # L85: 20834216 -> 20834660
$ws.Cells.Item(85, 13).Value = -57.71419999999989   # M85: -1863.5715 -> -57.71419999999989
$ws.Cells.Item(85, 14).Value = -20837156   # N85: -20836712 -> -20837156
# Row 122
$ws.Cells.Item(122, 8).Value = 1000004   # H122: 145172 -> 1000004
$ws.Cells.Item(122, 9).Value = 1000004   # I122: 252251 -> 1000004
$ws.Cells.Item(122, 10).Value = 0   # J122: 2400 -> 0
$ws.Cells.Item(122, 11).Value = 3000012   # K122: 756753 -> 3000012
$ws.Cells.Item(122, 12).Value = 0   # L122: 7200 -> 0
$ws.Cells.Item(122, 13).Value = -2997562   # M122: -754303 -> -2997562
$ws.Cells.Item(122, 14).ClearContents()   # N122: remove cell (was -12100)
# Row 123
$ws.Cells.Item(123, 8).Value = 29383.455   # H123: 28264.889 -> 29383.455
$ws.Cells.Item(123, 10).Value = 29383.455   # J123: 28264.889 -> 29383.455
$ws.Cells.Item(123, 12).Value = 29383.455   # L123: 28264.889 -> 29383.455
$ws.Cells.Item(123, 14).Value = -39183.455   # N123: -38064.889 -> -39183.455
# Row 126
$ws.Cells.Item(126, 8).Value = 2250   # H126: 2199.7058 -> 2250
$ws.Cells.Item(126, 9).Value = 2000   # I126: 1807.5 -> 2000
$ws.Cells.Item(126, 10).Value = 2300   # J126: 3141 -> 2300
$ws.Cells.Item(126, 11).Value = 6000   # K126: 5422.5 -> 6000
$ws.Cells.Item(126, 12).Value = 6900   # L126: 9423 -> 6900
$ws.Cells.Item(126, 13).Value = -3530   # M126: -2952.5 -> -3530
$ws.Cells.Item(126, 14).Value = -11840   # N126: -14363 -> -11840

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 1139.8723   # H132: 1204.4186 -> 1139.8723
$ws.Cells.Item(132, 9).Value = 852.94116   # I132: 900.5161000000001 -> 852.94116
$ws.Cells.Item(132, 10).Value = 1890.3077   # J132: 1989.5 -> 1890.3077
$ws.Cells.Item(132, 11).Value = 2558.82348   # K132: 2701.5483 -> 2558.82348
$ws.Cells.Item(132, 12).Value = 5670.9231   # L132: 5968.5 -> 5670.9231
$ws.Cells.Item(132, 13).Value = -28.82348000000002   # M132: -171.5483000000004 -> -28.82348000000002
$ws.Cells.Item(132, 14).Value = -10730.9231   # N132: -11028.5 -> -10730.9231

